$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) store plain text values (e.g. "60.740.88",
# "  -3.53%  "); force text format first so Excel does not coerce simple
# decimal-looking strings (e.g. "80.30") into numbers and drop trailing zeros.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.740.88"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.908.73"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.54"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.16"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -6.39%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.908.58"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.75"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.73%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.98%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.51%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.41"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -6.32%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.389.72"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.683.83"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.73"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.907.33"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "431.22"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.08%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.14%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.34%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.30"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.90"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.42%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.68%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.31%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.61"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.44%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.15"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.13%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.32%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0862"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.15%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.20%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.34%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.71"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.64"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.293"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.22"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "374.65"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.80%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.667.21"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.88"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.29"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.107"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.11%  "
